# Add Darklords to excel sheet
#
# SPDS-JP ("Structure Deck: Rise of the Shadow Lords" JP id list) gains three
# new cards that were missing from the sheet:
#   - Darklord Mastema        -> goes on the existing row for id 100405032
#   - Altar of the Darklords  -> goes on the existing row for id 100405035
#   - Darklords' Temptation   -> starts a brand new block of four rows
#     (100405037-100405040), mirroring the ":"/"," marker columns used by
#     every other card row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPDS-JP")
$ws.Activate()

# Existing rows that only need a card name filled into column A.
$ws.Range("A33").Value = "Darklord Mastema"
$ws.Range("A36").Value = "Altar of the Darklords"

# New card block starting at row 38 (ids 100405037 - 100405040).
$ws.Range("A38").Value = "Darklords' Temptation"
$ws.Range("B38").Value = 100405037
$ws.Range("C38").Value = ":"
$ws.Range("E38").Value = ","

$ws.Range("B39").Value = 100405038
$ws.Range("C39").Value = ":"
$ws.Range("E39").Value = ","

$ws.Range("B40").Value = 100405039
$ws.Range("C40").Value = ":"
$ws.Range("E40").Value = ","

$ws.Range("B41").Value = 100405040
$ws.Range("C41").Value = ":"
$ws.Range("E41").Value = ","

# Match the author's updated view/selection (scrolled down to the new rows).
$win = $excel.ActiveWindow
$win.ScrollRow = 38
$win.ScrollColumn = 1
$ws.Range("F45").Select()
